$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value2 = '26.491.10'
$ws.Range('E2').Value2 = '  +1.60%  '
$ws.Range('D3').Value2 = '1.670.85'
$ws.Range('E3').Value2 = '  +1.44%  '
$ws.Range('E4').Value2 = '  +0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value2 = '220.28'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value2 = '  +2.09%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value2 = '0.5280'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value2 = '  +0.85%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value2 = '1.002'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value2 = '  +0.10%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value2 = '0.2678'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value2 = '  +2.40%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value2 = '0.06377'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value2 = '  +0.12%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value2 = '21.77'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value2 = '  +4.23%  '
$ws.Range('E11').Value2 = '  +1.72%  '
$ws.Range('B12').Value2 = 'Polkadot'
$ws.Range('C12').Value2 = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value2 = '4.488'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value2 = '  +1.34%  '
$ws.Range('B13').Value2 = 'WrappedEther'
$ws.Range('C13').Value2 = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value2 = '1.668.99'
$ws.Range('E13').Value2 = '  +1.32%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value2 = '0.5566'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value2 = '  +0.27%  '
$ws.Range('D15').Value2 = '0.0₅8278'
$ws.Range('E15').Value2 = '  -0.61%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value2 = '65.60'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value2 = '  +0.97%  '
$ws.Range('D17').Value2 = '26.500.41'
$ws.Range('E17').Value2 = '  +1.64%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value2 = '4.761'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value2 = '  +0.58%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value2 = '193.06'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value2 = '  +2.38%  '
$ws.Range('E21').Value2 = '  +1.09%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value2 = '6.307'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value2 = '  +0.48%  '
$ws.Range('E23').Value2 = '  +0.15%  '
$ws.Range('E24').Value2 = '  +3.88%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value2 = '138.35'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value2 = '  -5.26%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value2 = '7.393'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value2 = '16.31'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value2 = '  +2.91%  '
$ws.Range('E28').Value2 = '  +2.10%  '
$ws.Range('E29').Value2 = '  +4.52%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value2 = '1.286'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value2 = '  +1.31%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value2 = '3.610'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value2 = '  +5.80%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value2 = '3.417'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value2 = '  +0.40%  '
$ws.Range('E33').Value2 = '  +1.78%  '
$ws.Range('E34').Value2 = '  +0.84%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value2 = '0.6172'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value2 = '  +9.41%  '
$ws.Range('E36').Value2 = '  +1.26%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value2 = '2.784'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value2 = '  +1.11%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value2 = '0.01615'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value2 = '  +0.33%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value2 = '6.039'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value2 = '  +3.40%  '
$ws.Range('D40').Value2 = '1.092.85'
$ws.Range('E40').Value2 = '  +6.13%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value2 = '0.8574'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value2 = '  -0.08%  '
$ws.Range('E42').Value2 = '  +0.04%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value2 = '100.64'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value2 = '  +1.90%  '
$ws.Range('D44').Value2 = '1.815.57'
$ws.Range('E44').Value2 = '  +1.06%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value2 = '58.56'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value2 = '  +4.82%  '
$ws.Range('B46').Value2 = 'EnergySwap'
$ws.Range('C46').Value2 = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value2 = '8.171'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value2 = '  +1.14%  '
$ws.Range('B47').Value2 = 'BabyDogeCoin'
$ws.Range('C47').Value2 = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D47').Value2 = '0.0₈104'
$ws.Range('E47').Value2 = '  -6.42%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value2 = '0.9982'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value2 = '  -0.63%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value2 = '1.515'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value2 = '  +9.41%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value2 = '0.05190'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value2 = '  +0.73%  '
$ws.Range('E51').Value2 = '  +0.44%  '
